$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51:58 down to 52:59.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the latest weekly price record.
$ws.Cells.Item(51,1).Value = 10
$ws.Cells.Item(51,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(51,3).Value = "La Araucanía"
$ws.Cells.Item(51,4).Value = 44637
$ws.Cells.Item(51,5).Value = 9
$ws.Cells.Item(51,6).Value = 300000001
$ws.Cells.Item(51,7).Value = "Rabanito"
$ws.Cells.Item(51,8).Value = "Sin especificar"
$ws.Cells.Item(51,9).Value = "Primera"
$ws.Cells.Item(51,10).Value = 30
$ws.Cells.Item(51,11).Value = 7000
$ws.Cells.Item(51,12).Value = 7000
$ws.Cells.Item(51,13).Value = 7000
$ws.Cells.Item(51,14).Value = "$/docena de paquetes"
$ws.Cells.Item(51,15).Value = "Provincia de Cautín"
$ws.Cells.Item(51,16).Value = 583
$ws.Cells.Item(51,17).Value = 12
$ws.Cells.Item(51,18).Value = "Hortaliza"
